$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 38: new diary entry (Programmieren, 1 hour, 2019-08-17) ---
$ws.Range("E37").Copy()
$ws.Range("E38").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E38").Value = 43694
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = "Stunden"
$ws.Range("H38").Value = "Programmieren"

# --- Row 39: new diary entry (Programmieren, 5 hours, 2019-08-18) ---
$ws.Range("E37").Copy()
$ws.Range("E39").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E39").Value = 43695
$ws.Range("F39").Value = 5
$ws.Range("G39").Value = "Stunden"
$ws.Range("H39").Value = "Programmieren"
$ws.Range("I39").Value = "Design neuer Einstellungen Klasse"

$excel.CutCopyMode = $false

# Update the selection/active cell to reflect the new last entry
$ws.Range("I39").Select()

$wb.Save()
